$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Edit existing staff member ALANA: last name TEST -> DANSKIN, role STAFF -> MANAGER
$ws.Range("B3").Value = "DANSKIN"
$ws.Range("D3").Value = "MANAGER"

# Widen the ROLE column (D) to fit the new, longer role text
$ws.Columns.Item(4).ColumnWidth = 9.333333333333334

Write-Host "Staff sheet updated"
